$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 4 new rows right before the old row 108 (pushes old 108:111 -> 112:115) ---
$ws.Rows.Item(108).Resize(4).Insert()

# The insert copies formatting (incl. the H column) from the row above (107) into the
# new blank rows. The new "oiltank" trailer rows don't use column H at all, so drop
# those cells completely (Clear removes both content + formatting -> no <c> emitted).
$ws.Range("H108:H111").Clear()

# --- Row 108: scout_trailer_oiltank ---
$ws.Cells.Item(108, 1).Value = 3
$ws.Cells.Item(108, 2).Value = "[media]\classes\trucks\trailers\scout_trailer_oiltank.xml"
$ws.Cells.Item(108, 3).Value = 900
$ws.Cells.Item(108, 4).Formula = "=IF(A108=3, C108*0.82,C108)"
$ws.Cells.Item(108, 6).Value = 900
$ws.Cells.Item(108, 7).Value = 900
$ws.Cells.Item(108, 9).Formula = "=G108/F108"

# --- Row 109: semitrailer_heavy_oiltank ---
$ws.Cells.Item(109, 1).Value = 3
$ws.Cells.Item(109, 2).Value = "[media]\classes\trucks\trailers\semitrailer_heavy_oiltank.xml"
$ws.Cells.Item(109, 3).Value = 42000
$ws.Cells.Item(109, 4).Formula = "=IF(A109=3, C109*0.82,C109)"
$ws.Cells.Item(109, 6).Value = 5000
$ws.Cells.Item(109, 7).Value = 5000
$ws.Cells.Item(109, 9).Formula = "=G109/F109"

# --- Row 110: semitrailer_oiltank ---
$ws.Cells.Item(110, 1).Value = 3
$ws.Cells.Item(110, 2).Value = "[media]\classes\trucks\trailers\semitrailer_oiltank.xml"
$ws.Cells.Item(110, 3).Value = 36000
$ws.Cells.Item(110, 4).Formula = "=IF(A110=3, C110*0.82,C110)"
$ws.Cells.Item(110, 6).Value = 3700
$ws.Cells.Item(110, 7).Value = 3700
$ws.Cells.Item(110, 9).Formula = "=G110/F110"

# --- Row 111: trailer_oiltank ---
$ws.Cells.Item(111, 1).Value = 3
$ws.Cells.Item(111, 2).Value = "[media]\classes\trucks\trailers\trailer_oiltank.xml"
$ws.Cells.Item(111, 3).Value = 20000
$ws.Cells.Item(111, 4).Formula = "=IF(A111=3, C111*0.82,C111)"
$ws.Cells.Item(111, 6).Value = 2000
$ws.Cells.Item(111, 7).Value = 3700
$ws.Cells.Item(111, 9).Formula = "=G111/F111"

# --- Stray formatted-but-empty cells left behind in row 2 (I2:J2:K2), matching the
# regular data-row style (picked up via a format-only paste from A3). ---
$ws.Range("A3").Copy()
$ws.Range("I2:K2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selection bookkeeping: the whole newly (re)entered block is highlighted ---
$ws.Range("A3:D115").Select()
